# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2 "H") with Wild Card round stats added in
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 263
$wsOff.Range("C2").Value = 177
$wsOff.Range("D2").Value = 61
$wsOff.Range("E2").Value = 35
$wsOff.Range("F2").Value = 4

# Update DEF sheet (row 2 "H") with Wild Card round stats added in
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 236
$wsDef.Range("C2").Value = 154
$wsDef.Range("D2").Value = 61
$wsDef.Range("E2").Value = 34
$wsDef.Range("F2").Value = 7

$wb.Save()
